$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new staff record (row 34) for the Lab Instructor
# added as part of the civil/ECE department roster update.
$ws.Range("A34").Value = "Mrs. RAJESWARI T"
$ws.Range("B34").Value = "Lab Instructor"
$ws.Range("C34").Value = "/static/images/profile_photos/009/VEC-009-05-006.webp"
$ws.Range("J34").Value = "VEC-009-05-006"

# Reflect the cursor/selection position left after the edit.
$ws.Range("H36").Select() | Out-Null
